$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "Heat Transfer" header to clarify it is q_gen - q_loss
$ws.Range("O1").Value = "Heat Transfer (q_gen - q_loss)"

# Add two new columns: Bulk Mean Temperature and Heat Transfer (BMT)
$ws.Range("P1").Value = "Bulk Mean Temperature"
$ws.Range("Q1").Value = "Heat Transfer (BMT)"

# Copy header style (bold, bordered, centered) from O1 onto the new headers
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

$ws.Range("P2").Value = 1.944074140543679
$ws.Range("Q2").Value = 1424.052505355127

$ws.Range("P3").Value = 5.020557813103682
$ws.Range("Q3").Value = 3665.178925903648

$ws.Range("P4").Value = 6.721114306937997
$ws.Range("Q4").Value = 3577.242197580631

$ws.Range("P5").Value = 2.920085621750038
$ws.Range("Q5").Value = 1562.348859539467

# Adjust column widths to accommodate the new, longer headers
$ws.Columns.Item(15).ColumnWidth = 30.7109375
$ws.Columns.Item(16).ColumnWidth = 21.7109375
$ws.Columns.Item(17).ColumnWidth = 19.7109375
